$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 301 - this shifts the existing
# rows 301-303 down to rows 303-305 (keeping their original values).
$ws.Rows.Item(301).Insert()
$ws.Rows.Item(301).Insert()

# Fill new row 301 with the new data point (date 2022-04-05 / Hayward / Primera)
$ws.Cells.Item(301, 1).Value = 8
$ws.Cells.Item(301, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(301, 3).Value = "Coquimbo"
$ws.Cells.Item(301, 4).Value = 44656
$ws.Cells.Item(301, 5).Value = 4
$ws.Cells.Item(301, 6).Value = "Fruta"
$ws.Cells.Item(301, 7).Value = 100101
$ws.Cells.Item(301, 8).Value = "Berries"
$ws.Cells.Item(301, 9).Value = 100101007
$ws.Cells.Item(301, 10).Value = "Kiwi"
$ws.Cells.Item(301, 11).Value = "Hayward"
$ws.Cells.Item(301, 12).Value = "Primera"
$ws.Cells.Item(301, 13).Value = 20
$ws.Cells.Item(301, 14).Value = 335000
$ws.Cells.Item(301, 15).Value = 340000
$ws.Cells.Item(301, 16).Value = 337500
$ws.Cells.Item(301, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(301, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(301, 19).Value = 750
$ws.Cells.Item(301, 20).Value = 450

# Fill new row 302 with the new data point (date 2022-04-05 / Hayward / Segunda)
$ws.Cells.Item(302, 1).Value = 8
$ws.Cells.Item(302, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(302, 3).Value = "Coquimbo"
$ws.Cells.Item(302, 4).Value = 44656
$ws.Cells.Item(302, 5).Value = 4
$ws.Cells.Item(302, 6).Value = "Fruta"
$ws.Cells.Item(302, 7).Value = 100101
$ws.Cells.Item(302, 8).Value = "Berries"
$ws.Cells.Item(302, 9).Value = 100101007
$ws.Cells.Item(302, 10).Value = "Kiwi"
$ws.Cells.Item(302, 11).Value = "Hayward"
$ws.Cells.Item(302, 12).Value = "Segunda"
$ws.Cells.Item(302, 13).Value = 14
$ws.Cells.Item(302, 14).Value = 285000
$ws.Cells.Item(302, 15).Value = 290000
$ws.Cells.Item(302, 16).Value = 287500
$ws.Cells.Item(302, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(302, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(302, 19).Value = 639
$ws.Cells.Item(302, 20).Value = 450
